# Atualização automática de preços de eletricidade
# Updates row 2 of the spot price table with the new day's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45918
$ws.Range("B2").Value = 109.18
$ws.Range("C2").Value = 105.26
$ws.Range("D2").Value = 104.54
$ws.Range("E2").Value = 104.78
$ws.Range("F2").Value = 104.78
$ws.Range("G2").Value = 105.26
$ws.Range("H2").Value = 107.28
$ws.Range("I2").Value = 126.55
$ws.Range("J2").Value = 111.4
$ws.Range("K2").Value = 104.26
$ws.Range("L2").Value = 56.4
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 50.08
$ws.Range("O2").Value = 50
$ws.Range("P2").Value = 44.02
$ws.Range("Q2").Value = 44
$ws.Range("R2").Value = 55
$ws.Range("S2").Value = 66.3
$ws.Range("T2").Value = 96.86
$ws.Range("U2").Value = 124.7
$ws.Range("V2").Value = 160
$ws.Range("W2").Value = 200
$ws.Range("X2").Value = 130
$ws.Range("Y2").Value = 107.65
$ws.Range("Z2").Value = 97.01
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 149.41
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 180
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 118.82
$ws.Range("AG2").Value = "10h-18h"
